# New weekly price observation: insert a row at row 17 (pushing the
# existing rows 17-43 down to 18-44) and populate it with the latest
# "Achicoria" reading for Vega Central Mapocho de Santiago.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("17:17").Insert()

$ws.Cells(17,1).Value2  = 9
$ws.Cells(17,2).Value2  = "Vega Central Mapocho de Santiago"
$ws.Cells(17,3).Value2  = "Metropolitana"
$ws.Cells(17,4).Value2  = 44994
$ws.Cells(17,5).Value2  = 13
$ws.Cells(17,6).Value2  = 100112010
$ws.Cells(17,7).Value2  = "Achicoria"
$ws.Cells(17,8).Value2  = "Sin especificar"
$ws.Cells(17,9).Value2  = "Primera"
$ws.Cells(17,10).Value2 = 70
$ws.Cells(17,11).Value2 = 7000
$ws.Cells(17,12).Value2 = 7000
$ws.Cells(17,13).Value2 = 7000
$ws.Cells(17,14).Value2 = "$/caja 16 unidades"
$ws.Cells(17,15).Value2 = "Provincia de Quillota"
$ws.Cells(17,16).Value2 = 438
$ws.Cells(17,17).Value2 = 16
$ws.Cells(17,18).Value2 = "Hortaliza"
